$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 21; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Text -eq "Consumer Discretionary") {
        $cell.Value = "Consumer Cyclical"
    }
}
